$wb = $excel.ActiveWorkbook

# --- Sheet "r Workspace_BusinessLine" (3rd tab) ---------------------------
# Order of writes matters: it determines the order new strings are appended
# to the shared-string table, which must match the target file.
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = "BL_BE_000001"
$ws3.Range("F4").Value = "BL_BE_000004"
$ws3.Range("B4").Value = "WS_BIB_BL_BE_000004"
$ws3.Range("C4").Value = "WS_BIB_BL_BE_000004"
$ws3.Range("B3").Value = "WS_BIB_BL_BE_000001"
$ws3.Range("C3").Value = "WS_BIB_BL_BE_000001"
$ws3.Range("E3").Value = "WS_BIB_CORPORATE"
$ws3.Range("E4").Value = "WS_BIB_RETAIL"
$ws3.Range("A3").Value = "CREATE/MODIFY"
$ws3.Range("A4").Value = "CREATE/MODIFY"

# --- Sheet "Workspace" (1st tab) ------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B3").Value = "WS_BIB_CORPORATE"
$ws1.Range("C3").Value = "Corporate"
$ws1.Range("D3").Value = "Corporate"
$ws1.Range("B4").Value = "WS_BIB_RETAIL"
$ws1.Range("C4").Value = "Retail"
$ws1.Range("D4").Value = "COUNTERPARTY_BIB"
$ws1.Range("A4").Value = "CREATE/MODIFY"
$ws1.Range("E4").Value = 17

# --- Sheet "r Workspace_TargetVariable" (4th tab) -------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B3").Value = "WS_BIB_CORPORATE_BE_TARGET"
$ws4.Range("C3").Value = "WS_BIB_CORPORATE_BE_TARGET"
$ws4.Range("D3").Value = "WS_BIB_CORPORATE_BE_TARGET"
$ws4.Range("E3").Value = "WS_BIB_CORPORATE"
$ws4.Range("B4").Value = "WS_BIB_RETAIL_BE_TARGET"
$ws4.Range("C4").Value = "WS_BIB_RETAIL_BE_TARGET"
$ws4.Range("D4").Value = "WS_BIB_RETAIL_BE_TARGET"
$ws4.Range("E4").Value = "WS_BIB_RETAIL"
$ws4.Range("A4").Value = "CREATE/MODIFY"
# F3 already carries the bold/red "framework" style (cellXfs index 3) -
# copy that formatting onto F4 before writing its value.
$ws4.Range("F3").Copy()
$ws4.Range("F4").PasteSpecial(-4122)
$ws4.Range("F4").Value = "BE_TARGET"

# --- Sheet "r Workspace_AnalysisUnit" (2nd tab) ---------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B3").Value = "WS_BIB_CORPORATE_COUNTERPARTY_BIB"
$ws2.Range("C3").Value = "WS_BIB_CORPORATE_COUNTERPARTY_BIB"
$ws2.Range("D3").Value = "WS_BIB_CORPORATE_COUNTERPARTY_BIB"
$ws2.Range("E3").Value = "WS_BIB_CORPORATE"
$ws2.Range("B4").Value = "WS_BIB_RETAIL_COUNTERPARTY_BIB"
$ws2.Range("C4").Value = "WS_BIB_RETAIL_COUNTERPARTY_BIB"
$ws2.Range("D4").Value = "WS_BIB_RETAIL_COUNTERPARTY_BIB"
$ws2.Range("E4").Value = "WS_BIB_RETAIL"
$ws2.Range("A4").Value = "CREATE/MODIFY"
$ws2.Range("F4").Value = "COUNTERPARTY_BIB"

# --- Column width adjustments ---------------------------------------------
# (ColumnWidth is expressed in "characters"; the engine snaps to its own
# internal pixel grid, so these are the closest achievable values to the
# target OOXML widths.)
$ws1.Columns.Item(2).ColumnWidth = 18.833333333333332

$ws2.Columns.Item(2).ColumnWidth = 38.5
$ws2.Columns.Item(3).ColumnWidth = 38.5
$ws2.Columns.Item(4).ColumnWidth = 38.5
$ws2.Columns.Item(5).ColumnWidth = 18.833333333333332

$ws3.Columns.Item(2).ColumnWidth = 20.333333333333332

$ws4.Columns.Item(2).ColumnWidth = 30.166666666666668
$ws4.Columns.Item(3).ColumnWidth = 30.166666666666668
$ws4.Columns.Item(4).ColumnWidth = 30.166666666666668
$ws4.Columns.Item(5).ColumnWidth = 18.833333333333332

# --- Selections / active cells (set last, in tab order, sheet4 last so it
#     remains the active tab as in the source file) ------------------------
$ws1.Range("B3:B4").Select()
$ws2.Range("D10").Select()
$ws3.Range("E7").Select()
$ws4.Range("D7").Select()
